$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column U (2023) mirrors the existing 2022 (column T) formatting exactly:
# copy T4:T14 formats into U4:U14 first, then write the new values on top.
$ws.Range("T4:T14").Copy()
$ws.Range("U4:U14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("U4").Value = 2023
$ws.Range("U5").Value = 0.5
$ws.Range("U6").Value = 0.3
$ws.Range("U7").Value = 0.4
$ws.Range("U8").Value = 0.4
$ws.Range("U9").Value = 3.2
$ws.Range("U10").Value = 0.6
$ws.Range("U11").Value = "-"
$ws.Range("U12").Value = 0.6
$ws.Range("U13").Value = 0.1
$ws.Range("U14").Value = 0.5

# Match the saved selection state (B1 instead of the prior V7).
$ws.Range("B1").Select() | Out-Null
